$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 18 data, following the pattern of the preceding rows
$ws.Range("A18").Value = "Exp 22"
$ws.Range("B18").Value = 0.3
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = "Local"
$ws.Range("E18").Value = -1
$ws.Range("F18").Value = "Exp 22.png"

# Match the style of the row above (A17:E17 -> A18:E18)
$ws.Range("A17:E17").Copy()
$ws.Range("A18:E18").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Match the row height used by the preceding rows
$ws.Rows.Item(18).RowHeight = $ws.Rows.Item(17).RowHeight

# Update the selection to reflect the active cell shown in the diff
$ws.Range("B18").Select()
